$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 updates
$ws.Range("G5").Value = 1.85
$ws.Range("H5").Value = 3.1
$ws.Range("I5").Value = 4.45
$ws.Range("M5").Value = 2.35
$ws.Range("N5").Value = 2.32
$ws.Range("O5").Value = 1.47
$ws.Range("P5").Value = 1.5
$ws.Range("Q5").Value = 2.25
$ws.Range("R5").Value = 2.1
$ws.Range("S5").Value = 1.57
$ws.Range("T5").Value = 5.3
$ws.Range("U5").Value = 7.5
$ws.Range("V5").Value = 8.75
$ws.Range("W5").Value = 15
$ws.Range("Y5").Value = 40
$ws.Range("Z5").Value = 6.6
$ws.Range("AA5").Value = 6.2
$ws.Range("AB5").Value = 20
$ws.Range("AC5").Value = 120
$ws.Range("AF5").Value = 23
$ws.Range("AG5").Value = 16
$ws.Range("AH5").Value = 80
$ws.Range("AI5").Value = 60
$ws.Range("AJ5").Value = 75

# Row 6 updates
$ws.Range("G6").Value = 2.22
$ws.Range("H6").Value = 2.75
$ws.Range("I6").Value = 3.65
$ws.Range("L6").Value = 1.52
$ws.Range("M6").Value = 2.2
$ws.Range("N6").Value = 2.5
$ws.Range("O6").Value = 1.4
$ws.Range("P6").Value = 1.55
$ws.Range("Q6").Value = 2.15
$ws.Range("R6").Value = 2.05
$ws.Range("S6").Value = 1.6
$ws.Range("T6").Value = 5.5
$ws.Range("U6").Value = 9.25
$ws.Range("V6").Value = 9.5
$ws.Range("W6").Value = 22
$ws.Range("X6").Value = 23
$ws.Range("Y6").Value = 45
$ws.Range("Z6").Value = 5.8
$ws.Range("AA6").Value = 5.6
$ws.Range("AB6").Value = 18
$ws.Range("AC6").Value = 120
$ws.Range("AF6").Value = 18
$ws.Range("AG6").Value = 13
$ws.Range("AH6").Value = 60
$ws.Range("AI6").Value = 45
$ws.Range("AJ6").Value = 60
